{"js": "// Adds a new \"meta description\" block (bold title + italic description)\n// at the end of the document, and removes the old one that used to sit\n// right under the top heading (it contained \"Meta description: ...\").\n\nconst body = context.document.body;\n\n// --- Step 1: remove the old \"Meta description\" paragraph. Locate it by\n// its distinctive text rather than a fixed index, so the script is\n// resilient to the exact paragraph layout.\nconst metaSearchResults = body.search(\"Meta description\", { matchCase: false });\nmetaSearchResults.load(\"text\");\nawait context.sync();\n\nif (metaSearchResults.items.length > 0) {\n  const metaDescriptionParagraph = metaSearchResults.items[0].paragraphs.getFirst();\n  metaDescriptionParagraph.delete();\n  await context.sync();\n}\n\n// --- Step 2: replace the final paragraph (the italic AI image-prompt\n// paragraph, \"Create a cartoon-style feature image...\") with two new\n// paragraphs: a bold title line and an italic meta-description line.\n// Locate it by its distinctive text; fall back to the very last\n// paragraph in the body if, for some reason, it can't be found.\nconst promptSearchResults = body.search(\"Create a cartoon-style feature image\", { matchCase: false });\npromptSearchResults.load(\"text\");\nawait context.sync();\n\nlet oldImagePromptParagraph;\nif (promptSearchResults.items.length > 0) {\n  oldImagePromptParagraph = promptSearchResults.items[0].paragraphs.getFirst();\n} else {\n  const paragraphs = body.paragraphs;\n  paragraphs.load(\"text\");\n  await context.sync();\n  oldImagePromptParagraph = paragraphs.items[paragraphs.items.length - 1];\n}\n\nconst newContentHtml =\n  \"<p><b>Play Fruit Party slot game for free</b></p>\" +\n  \"<p><i>Discover the features and potential winnings of Fruit Party slot game. Play for free and experience its cascade system and medium volatility.</i></p>\";\n\noldImagePromptParagraph.insertHtml(newContentHtml, Word.InsertLocation.before);\nawait context.sync();\n\n// The original last paragraph (image-prompt text) got pushed down to the\n// very end again by the insert above \u2014 re-fetch and delete it.\nconst refreshedParagraphs = body.paragraphs;\nrefreshedParagraphs.load(\"text\");\nawait context.sync();\n\nconst trailingParagraph = refreshedParagraphs.items[refreshedParagraphs.items.length - 1];\ntrailingParagraph.delete();\nawait context.sync();\n", "ps1": "# Adds a new \"meta description\" block (bold title + italic description)\n# at the end of the document, and removes the old one that used to sit\n# right under the top heading (it contained \"Meta description: ...\").\n\n$d = $word.ActiveDocument\n\nfunction Get-ParagraphIndexContaining {\n    param([string]$needle)\n    $count = $d.Paragraphs.Count\n    for ($i = 1; $i -le $count; $i++) {\n        $p = $d.Paragraphs.Item($i)\n        if ($p.Range.Text -like \"*$needle*\") {\n            return $i\n        }\n    }\n    return -1\n}\n\n# --- Step 1: remove the old \"Meta description\" paragraph. Locate it by\n# its distinctive text rather than assuming a fixed paragraph index.\n$metaDescriptionIndex = Get-ParagraphIndexContaining(\"Meta description\")\nif ($metaDescriptionIndex -ge 1) {\n    $metaDescriptionParagraph = $d.Paragraphs.Item($metaDescriptionIndex)\n    $metaDescriptionParagraph.Range.Delete()\n}\n\n# --- Step 2: append two new paragraphs where the old trailing AI\n# image-prompt paragraph (\"Create a cartoon-style feature image...\")\n# used to be: a bold title line and an italic meta-description line.\n$imagePromptIndex = Get-ParagraphIndexContaining(\"Create a cartoon-style feature image\")\nif ($imagePromptIndex -lt 1) {\n    # Fallback: if the text can't be found for some reason, target the\n    # very last paragraph of the document.\n    $imagePromptIndex = $d.Paragraphs.Count\n}\n\n# Anchor the insertion on the paragraph right before the image-prompt\n# one, so the new paragraphs don't inherit the italic run formatting\n# that lives on the image-prompt paragraph itself.\n$anchorParagraph = $d.Paragraphs.Item($imagePromptIndex - 1)\n\n$anchorRange = $anchorParagraph.Range\n$anchorRange.Collapse(0)   # wdCollapseEnd\n$anchorRange.InsertParagraphAfter()\n$titleParaIndex = $imagePromptIndex\n\n$titlePara = $d.Paragraphs.Item($titleParaIndex)\n$titleRange = $titlePara.Range\n$titleRange.Collapse(0)    # wdCollapseEnd\n$titleRange.InsertParagraphAfter()\n$descriptionParaIndex = $titleParaIndex + 1\n\n# Normalize both freshly inserted (still empty) paragraphs to the\n# \"Normal\" style so they don't keep any inherited list/heading style.\n$titlePara2 = $d.Paragraphs.Item($titleParaIndex)\n$titlePara2.Style = \"Normal\"\n$descriptionPara2 = $d.Paragraphs.Item($descriptionParaIndex)\n$descriptionPara2.Style = \"Normal\"\n\n# Fill in the bold title paragraph (text excludes the trailing\n# paragraph mark so the mark itself stays unformatted).\n$titlePara3 = $d.Paragraphs.Item($titleParaIndex)\n$titleTextRange = $d.Range($titlePara3.Range.Start, $titlePara3.Range.End - 1)\n$titleTextRange.Text = \"Play Fruit Party slot game for free\"\n$titleTextRange2 = $d.Range($titlePara3.Range.Start, $titlePara3.Range.End - 1)\n$titleTextRange2.Font.Bold = 1\n\n# Fill in the italic description paragraph.\n$descriptionPara3 = $d.Paragraphs.Item($descriptionParaIndex)\n$descriptionTextRange = $d.Range($descriptionPara3.Range.Start, $descriptionPara3.Range.End - 1)\n$descriptionTextRange.Text = \"Discover the features and potential winnings of Fruit Party slot game. Play for free and experience its cascade system and medium volatility.\"\n$descriptionTextRange2 = $d.Range($descriptionPara3.Range.Start, $descriptionPara3.Range.End - 1)\n$descriptionTextRange2.Font.Italic = 1\n\n# --- Step 3: remove the old trailing paragraph (the AI image-prompt\n# text), which was pushed down by the two inserts above.\n$oldImagePromptIndex = $descriptionParaIndex + 1\n$oldImagePromptParagraph = $d.Paragraphs.Item($oldImagePromptIndex)\n$oldImagePromptParagraph.Range.Delete()\n"}
